# Add new "BS" currency rows (rows 14-25) mirroring the existing monthly
# dates already present in rows 2-13, with avg_closed value of 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date serials (1900 date system) for the 1st of each month, Jan-Dec 2025,
# matching the values already used for the USD rows.
$dates = @(45658, 45689, 45717, 45748, 45778, 45809, 45839, 45870, 45901, 45931, 45962, 45992)

# Copy the date cell formatting (numFmtId 14) from an existing row so we
# reuse the same style index instead of creating a brand new one.
$ws.Cells.Item(2, 1).Copy() | Out-Null

$row = 14
foreach ($d in $dates) {
    $ws.Cells.Item($row, 1).Value = $d
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($row, 2).Value = "BS"
    $ws.Cells.Item($row, 3).Value = 1
    $row = $row + 1
}
$excel.CutCopyMode = $false

# Update selection to match the newly entered range.
$ws.Range("C14:C25").Select()
